$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 463; this shifts the existing rows 463:556 down to
# 464:557 (and carries their formatting along), matching the diff's "new
# record pushed the weekly series down by one row" shape.
$ws.Rows("463:463").Insert()

# Populate the newly inserted row 463 with the new weekly record.
$ws.Cells.Item(463, 1).Value  = 3
$ws.Cells.Item(463, 2).Value  = 'Femacal de La Calera'
$ws.Cells.Item(463, 3).Value  = 'Coquimbo'
$ws.Cells.Item(463, 4).Value  = 45258
$ws.Cells.Item(463, 5).Value  = 5
$ws.Cells.Item(463, 6).Value  = 100112001
$ws.Cells.Item(463, 7).Value  = 'Berenjena'
$ws.Cells.Item(463, 8).Value  = 'Sin especificar'
$ws.Cells.Item(463, 9).Value  = 'Primera'
$ws.Cells.Item(463, 10).Value = 75
$ws.Cells.Item(463, 11).Value = 10000
$ws.Cells.Item(463, 12).Value = 11000
$ws.Cells.Item(463, 13).Value = 10467
$ws.Cells.Item(463, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(463, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(463, 16).Value = 174
$ws.Cells.Item(463, 17).Value = 60
$ws.Cells.Item(463, 18).Value = 'Hortaliza'
